$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The published content got reshuffled: the long Objetivos/Programa/Bibliografia
# paragraphs were dropped and the remaining labels + short values slid up to fill
# rows 13-23, which is three rows shorter than before (26 -> 23).

# Row 10
$ws.Cells.Item(10,"B").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Cells.Item(10,"C").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# Row 13
$ws.Cells.Item(13,"A").Value = "Programa resumido:"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Cells.Item(14,"A").Value = "Short syllabus:"
$ws.Cells.Item(14,"B").Clear()
$ws.Cells.Item(14,"C").Clear()
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Cells.Item(15,"A").Value = "Programa:"
$ws.Cells.Item(15,"B").Value = "3480026 - João Paulo Pascon"
$ws.Cells.Item(15,"C").Value = "3480026 - João Paulo Pascon"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Cells.Item(16,"A").Value = "Syllabus:"
$ws.Cells.Item(16,"B").Clear()
$ws.Cells.Item(16,"C").Clear()
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Cells.Item(17,"A").Value = "Avaliação:"
$ws.Rows.Item(17).EntireRow.AutoFit()

# Row 18
$ws.Cells.Item(18,"A").Value = "Método:"
$ws.Cells.Item(18,"B").Value = "7797767 - Viktor Pastoukhov"
$ws.Cells.Item(18,"C").Value = "7797767 - Viktor Pastoukhov"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Cells.Item(19,"A").Value = "Critério:"
$ws.Cells.Item(19,"B").Value = "A avaliação será composta por duas provas (P1 e P2)."
$ws.Cells.Item(10,"B").Copy()
$ws.Cells.Item(19,"B").PasteSpecial(-4122)
$ws.Cells.Item(19,"C").Value = "A avaliação será composta por duas provas (P1 e P2)."
$ws.Cells.Item(10,"C").Copy()
$ws.Cells.Item(19,"C").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Cells.Item(20,"A").Value = "Norma de recuperação:"
$ws.Cells.Item(20,"B").Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Cells.Item(10,"B").Copy()
$ws.Cells.Item(20,"B").PasteSpecial(-4122)
$ws.Cells.Item(20,"C").Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Cells.Item(10,"C").Copy()
$ws.Cells.Item(20,"C").PasteSpecial(-4122)
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Cells.Item(21,"A").Value = "Bibliografia:"
$ws.Cells.Item(21,"B").Value = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
$ws.Cells.Item(21,"C").Value = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Cells.Item(22,"A").Value = "Requisitos:"
$ws.Cells.Item(22,"B").Clear()
$ws.Cells.Item(22,"C").Clear()
$ws.Rows.Item(22).EntireRow.AutoFit()

# Row 23
$ws.Cells.Item(23,"A").Clear()
$ws.Cells.Item(23,"B").Value = "LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)`n"
$ws.Cells.Item(23,"C").Value = "LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# The old tail (rows 24-26) no longer exists now that its content folded into 13-23 above.
$ws.Range("A24:A26").EntireRow.Delete()

$excel.CutCopyMode = 0
